$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.320058704838459
$ws.Range("C2").Value = 0.3004919998410571
$ws.Range("D2").Value = 0.6481261315988718
$ws.Range("E2").Value = 0.2646925025571107
$ws.Range("G2").Value = 0.5743891366149256
$ws.Range("H2").Value = 0.7011943039804009
$ws.Range("J2").Value = 0.138051621387568
$ws.Range("O2").Value = 2.528277394862727
$ws.Range("B3").Value = 1.176917813548187
$ws.Range("C3").Value = 0.262335176208353
$ws.Range("D3").Value = 0.6374256006386076
$ws.Range("E3").Value = 0.2594462408059854
$ws.Range("G3").Value = 0.5786935525712096
$ws.Range("H3").Value = 0.7087774914066571
$ws.Range("J3").Value = 0.1345776010357795
$ws.Range("O3").Value = 2.552918075980159
$ws.Range("B4").Value = 1.088880517327368
$ws.Range("C4").Value = 0.238812927187297
$ws.Range("D4").Value = 0.6311935226437413
$ws.Range("E4").Value = 0.2563693357720567
$ws.Range("G4").Value = 0.5819458169097018
$ws.Range("H4").Value = 0.7139030765612517
$ws.Range("J4").Value = 0.1325241615265185
$ws.Range("O4").Value = 2.570309741029774
$ws.Range("B5").Value = 1.052969376192152
$ws.Range("C5").Value = 0.2292044505356614
$ws.Range("D5").Value = 0.6287389494842159
$ws.Range("E5").Value = 0.2551517544532445
$ws.Range("G5").Value = 0.5834237930548554
$ws.Range("H5").Value = 0.7161096957901592
$ws.Range("J5").Value = 0.1317073660583503
$ws.Range("O5").Value = 2.577964295056617
$ws.Range("B6").Value = 1.047004291507847
$ws.Range("C6").Value = 0.2276076014078967
$ws.Range("D6").Value = 0.6283365072840184
$ws.Range("E6").Value = 0.2549517671609109
$ws.Range("G6").Value = 0.5836784129486929
$ws.Range("H6").Value = 0.716483219745939
$ws.Range("J6").Value = 0.1315729449220271
$ws.Range("O6").Value = 2.579269546979759
$ws.Range("B7").Value = 1.088396347031903
$ws.Range("C7").Value = 0.2386834360842727
$ws.Range("D7").Value = 0.6311600750040043
$ws.Range("E7").Value = 0.2563527681479059
$ws.Range("G7").Value = 0.5819651320899482
$ws.Range("H7").Value = 0.7139323586257049
$ws.Range("J7").Value = 0.1325130649976884
$ws.Range("O7").Value = 2.570410678057257
$ws.Range("B8").Value = 1.270735658752301
$ws.Range("C8").Value = 0.2873553535975191
$ws.Range("D8").Value = 0.6443664322873133
$ws.Range("E8").Value = 0.2628536302567142
$ws.Range("G8").Value = 0.5757464879404353
$ws.Range("H8").Value = 0.7037114554139663
$ws.Range("J8").Value = 0.1368372374033342
$ws.Range("O8").Value = 2.536303059503609
$ws.Range("B9").Value = 1.627055317918007
$ws.Range("C9").Value = 0.3820341702927976
$ws.Range("D9").Value = 0.6729474358068899
$ws.Range("E9").Value = 0.2767483336182792
$ws.Range("G9").Value = 0.5684127272811139
$ws.Range("H9").Value = 0.687400356347041
$ws.Range("J9").Value = 0.145950416553525
$ws.Range("O9").Value = 2.48743802183526
$ws.Range("B10").Value = 1.888014564201171
$ws.Range("C10").Value = 0.451104359689964
$ws.Range("D10").Value = 0.6955861255478339
$ws.Range("E10").Value = 0.2876591857765121
$ws.Range("G10").Value = 0.5660248826096961
$ws.Range("H10").Value = 0.6777015531122146
$ws.Range("J10").Value = 0.1530353512828526
$ws.Range("O10").Value = 2.462621520224701
$ws.Range("B11").Value = 2.006539292993125
$ws.Range("C11").Value = 0.4824154037028165
$ws.Range("D11").Value = 0.7062423817162369
$ws.Range("E11").Value = 0.2927761847345138
$ws.Range("G11").Value = 0.565598032576375
$ws.Range("H11").Value = 0.6737876086447727
$ws.Range("J11").Value = 0.1563438134191983
$ws.Range("O11").Value = 2.453760236087703
$ws.Range("B12").Value = 2.051393055086066
$ws.Range("C12").Value = 0.4942558603582938
$ws.Range("D12").Value = 0.7103291102112053
$ws.Range("E12").Value = 0.2947359848193543
$ws.Range("G12").Value = 0.5655318480849161
$ws.Range("H12").Value = 0.672377301874036
$ws.Range("J12").Value = 0.1576089790116697
$ws.Range("O12").Value = 2.450755543825153
$ws.Range("B13").Value = 2.041734316891507
$ws.Range("C13").Value = 0.4917065436152939
$ws.Range("D13").Value = 0.709446672738693
$ws.Range("E13").Value = 0.2943129234714874
$ws.Range("G13").Value = 0.5655418466473634
$ws.Range("H13").Value = 0.6726778393998245
$ws.Range("J13").Value = 0.1573359544324262
$ws.Range("O13").Value = 2.45138702441281
$ws.Range("B14").Value = 2.010230031720141
$ws.Range("C14").Value = 0.4833898566017751
$ws.Range("D14").Value = 0.7065775685704807
$ws.Range("E14").Value = 0.2929369755803819
$ws.Range("G14").Value = 0.5655906708600611
$ws.Range("H14").Value = 0.6736701413828996
$ws.Range("J14").Value = 0.1564476522168263
$ws.Range("O14").Value = 2.453505995941867
$ws.Range("B15").Value = 1.990928911508263
$ws.Range("C15").Value = 0.4782934969366579
$ws.Range("D15").Value = 0.7048268574002918
$ws.Range("E15").Value = 0.2920970468135451
$ws.Range("G15").Value = 0.5656330268896284
$ws.Range("H15").Value = 0.6742873137338989
$ws.Range("J15").Value = 0.155905147131719
$ws.Range("O15").Value = 2.454849674798055
$ws.Range("B16").Value = 1.880264752639107
$ws.Range("C16").Value = 0.4490558548561694
$ws.Range("D16").Value = 0.6948969126462714
$ws.Range("E16").Value = 0.2873278700314259
$ws.Range("G16").Value = 0.5660661011879995
$ws.Range("H16").Value = 0.6779673780334861
$ws.Range("J16").Value = 0.1528208582837891
$ws.Range("O16").Value = 2.463249631305445
$ws.Range("B17").Value = 1.812326492997158
$ws.Range("C17").Value = 0.4310910517034472
$ws.Range("D17").Value = 0.6888968431438798
$ws.Range("E17").Value = 0.2844414861010733
$ws.Range("G17").Value = 0.5665011162663802
$ws.Range("H17").Value = 0.6803526892155247
$ws.Range("J17").Value = 0.1509506637975733
$ws.Range("O17").Value = 2.469025819868222
$ws.Range("B18").Value = 1.773232730717382
$ws.Range("C18").Value = 0.4207479186234195
$ws.Range("D18").Value = 0.6854794408133387
$ws.Range("E18").Value = 0.2827957707157225
$ws.Range("G18").Value = 0.5668133527012316
$ws.Range("H18").Value = 0.6817715280370322
$ws.Range("J18").Value = 0.1498830223158478
$ws.Range("O18").Value = 2.47257654796644
$ws.Range("B19").Value = 1.759993318084355
$ws.Range("C19").Value = 0.4172441699210481
$ws.Range("D19").Value = 0.6843281520071969
$ws.Range("E19").Value = 0.2822410421560662
$ws.Range("G19").Value = 0.5669297044899935
$ws.Range("H19").Value = 0.6822599666955824
$ws.Range("J19").Value = 0.1495229180859212
$ws.Range("O19").Value = 2.4738179407492
$ws.Range("B20").Value = 1.819560460321952
$ws.Range("C20").Value = 0.4330045014868915
$ws.Range("D20").Value = 0.6895320752620648
$ws.Range("E20").Value = 0.2847472503077029
$ws.Range("G20").Value = 0.5664483839773453
$ws.Range("H20").Value = 0.6800939163934459
$ws.Range("J20").Value = 0.1511489163630273
$ws.Range("O20").Value = 2.468387281649512
$ws.Range("B21").Value = 2.019484407047003
$ws.Range("C21").Value = 0.485833118089829
$ws.Range("D21").Value = 0.7074188981868303
$ws.Range("E21").Value = 0.2933405248086913
$ws.Range("G21").Value = 0.5655737343492291
$ws.Range("H21").Value = 0.6733767275373168
$ws.Range("J21").Value = 0.1567082334368877
$ws.Range("O21").Value = 2.452874065740247
$ws.Range("B22").Value = 2.149976390571794
$ws.Range("C22").Value = 0.5202640486709811
$ws.Range("D22").Value = 0.7194087703644527
$ws.Range("E22").Value = 0.2990855714625482
$ws.Range("G22").Value = 0.5655587750765676
$ws.Range("H22").Value = 0.6694053631666748
$ws.Range("J22").Value = 0.160413409090097
$ws.Range("O22").Value = 2.444781316817995
$ws.Range("B23").Value = 2.080346592228238
$ws.Range("C23").Value = 0.5018965699233036
$ws.Range("D23").Value = 0.7129821231791595
$ws.Range("E23").Value = 0.2960075369952122
$ws.Range("G23").Value = 0.5655156139336981
$ws.Range("H23").Value = 0.6714865790688265
$ws.Range("J23").Value = 0.1584293038580569
$ws.Range("O23").Value = 2.448912775965084
$ws.Range("B24").Value = 1.816290093033103
$ws.Range("C24").Value = 0.4321394772967437
$ws.Range("D24").Value = 0.6892447867793123
$ws.Range("E24").Value = 0.2846089716225961
$ws.Range("G24").Value = 0.566472030726743
$ws.Range("H24").Value = 0.6802107596625149
$ws.Range("J24").Value = 0.151059262838757
$ws.Range("O24").Value = 2.468675248707314
$ws.Range("B25").Value = 1.530802123718274
$ws.Range("C25").Value = 0.3565055272806035
$ws.Range("D25").Value = 0.6649277972228447
$ws.Range("E25").Value = 0.2728663491572334
$ws.Range("G25").Value = 0.569872486040822
$ws.Range("H25").Value = 0.6914123865399375
$ws.Range("J25").Value = 0.1434169020715146
$ws.Range("O25").Value = 2.498717892768212
